# Insert a new data row at row 223 (pushing the existing rows 223-322 down to 224-323),
# then populate the new row 223 with the new price-observation record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 223:322 down by one row, creating a blank row 223.
$ws.Rows("223:223").Insert()

# Populate the newly inserted row 223 with the new record
# (same "shape" as the surrounding Cilantro / Vega Modelo de Temuco rows).
$ws.Range("A223").Value = 10
$ws.Range("B223").Value = "Vega Modelo de Temuco"
$ws.Range("C223").Value = "La Araucanía"
$ws.Range("D223").Value = 44636
$ws.Range("E223").Value = 9
$ws.Range("F223").Value = 100112040
$ws.Range("G223").Value = "Cilantro"
$ws.Range("H223").Value = "Sin especificar"
$ws.Range("I223").Value = "Primera"
$ws.Range("J223").Value = 30
$ws.Range("K223").Value = 5000
$ws.Range("L223").Value = 5000
$ws.Range("M223").Value = 5000
$ws.Range("N223").Value = "$/docena de atados (2 kilos)"
$ws.Range("O223").Value = "Provincia de Cautín"
$ws.Range("P223").Value = 2500
$ws.Range("Q223").Value = 2
$ws.Range("R223").Value = "Hortaliza"
